# Update generated output numbers (and one row removal) per commit
# "Update gh-pages to output generated at 456a3b4"

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws3 = $wb.Worksheets.Item("本地生活")
$ws4 = $wb.Worksheets.Item("全部类型")

# -----------------------------------------------------------------
# Sheet: 展览 (Exhibition) - update "want to go" counts in column F
# -----------------------------------------------------------------
$ws1.Cells.Item(3, 6).Value = 145
$ws1.Cells.Item(4, 6).Value = 134
$ws1.Cells.Item(5, 6).Value = 1948
$ws1.Cells.Item(7, 6).Value = 4051
$ws1.Cells.Item(8, 6).Value = 527
$ws1.Cells.Item(9, 6).Value = 1040
$ws1.Cells.Item(12, 6).Value = 366
$ws1.Cells.Item(13, 6).Value = 96
$ws1.Cells.Item(14, 6).Value = 2174
$ws1.Cells.Item(16, 6).Value = 650156
$ws1.Cells.Item(17, 6).Value = 1603
$ws1.Cells.Item(18, 6).Value = 483
$ws1.Cells.Item(19, 6).Value = 1432
$ws1.Cells.Item(22, 6).Value = 1254
$ws1.Cells.Item(23, 6).Value = 2170
$ws1.Cells.Item(24, 6).Value = 1112
$ws1.Cells.Item(25, 6).Value = 2669
$ws1.Cells.Item(26, 6).Value = 1535
$ws1.Cells.Item(27, 6).Value = 771
$ws1.Cells.Item(28, 6).Value = 1513
$ws1.Cells.Item(29, 6).Value = 520
$ws1.Cells.Item(33, 6).Value = 2003
$ws1.Cells.Item(34, 6).Value = 1337
$ws1.Cells.Item(35, 6).Value = 1209
$ws1.Cells.Item(36, 6).Value = 2108
$ws1.Cells.Item(37, 6).Value = 1131
$ws1.Cells.Item(38, 6).Value = 39
$ws1.Cells.Item(39, 6).Value = 192
$ws1.Cells.Item(41, 6).Value = 2550

# -----------------------------------------------------------------
# Sheet: 演出 (Performance) - update "want to go" counts in column F
# -----------------------------------------------------------------
$ws2.Cells.Item(9, 6).Value = 95
$ws2.Cells.Item(10, 6).Value = 471
$ws2.Cells.Item(11, 6).Value = 144402
$ws2.Cells.Item(12, 6).Value = 144402
$ws2.Cells.Item(19, 6).Value = 331
$ws2.Cells.Item(21, 6).Value = 406
$ws2.Cells.Item(22, 6).Value = 406
$ws2.Cells.Item(23, 6).Value = 115
$ws2.Cells.Item(27, 6).Value = 526
$ws2.Cells.Item(28, 6).Value = 89
$ws2.Cells.Item(32, 6).Value = 321
$ws2.Cells.Item(33, 6).Value = 268

# -----------------------------------------------------------------
# Sheet: 本地生活 (Local Life) - row 5 ("「排球少年!!垃圾场决战 ×
# animate cafe」") was removed entirely; rows below it shift up by
# one, and their "want to go" counts (column F) were refreshed.
# -----------------------------------------------------------------
$ws3.Rows.Item(5).Delete() | Out-Null

$ws3.Cells.Item(4, 6).Value = 3116
$ws3.Cells.Item(5, 6).Value = 236
$ws3.Cells.Item(6, 6).Value = 13
$ws3.Cells.Item(7, 6).Value = 817
$ws3.Cells.Item(8, 6).Value = 1152
$ws3.Cells.Item(9, 6).Value = 631
$ws3.Cells.Item(10, 6).Value = 1577
$ws3.Cells.Item(11, 6).Value = 472
$ws3.Cells.Item(12, 6).Value = 59
$ws3.Cells.Item(13, 6).Value = 1824

# -----------------------------------------------------------------
# Sheet: 全部类型 (All Types) - update "want to go" counts in column F
# -----------------------------------------------------------------
$ws4.Cells.Item(2, 6).Value = 817
$ws4.Cells.Item(3, 6).Value = 631
$ws4.Cells.Item(4, 6).Value = 145
$ws4.Cells.Item(5, 6).Value = 1577
$ws4.Cells.Item(6, 6).Value = 472
$ws4.Cells.Item(7, 6).Value = 134
$ws4.Cells.Item(8, 6).Value = 1824
$ws4.Cells.Item(9, 6).Value = 4051
$ws4.Cells.Item(11, 6).Value = 527
$ws4.Cells.Item(14, 6).Value = 366
$ws4.Cells.Item(15, 6).Value = 2174
$ws4.Cells.Item(18, 6).Value = 650157
$ws4.Cells.Item(19, 6).Value = 95
$ws4.Cells.Item(20, 6).Value = 471
$ws4.Cells.Item(21, 6).Value = 1603
$ws4.Cells.Item(22, 6).Value = 144402
$ws4.Cells.Item(23, 6).Value = 483
$ws4.Cells.Item(24, 6).Value = 1432
$ws4.Cells.Item(27, 6).Value = 1254
$ws4.Cells.Item(28, 6).Value = 2170
$ws4.Cells.Item(29, 6).Value = 1112
$ws4.Cells.Item(30, 6).Value = 2669
$ws4.Cells.Item(31, 6).Value = 1535
$ws4.Cells.Item(33, 6).Value = 1513
$ws4.Cells.Item(34, 6).Value = 406
$ws4.Cells.Item(35, 6).Value = 520
$ws4.Cells.Item(36, 6).Value = 115
$ws4.Cells.Item(41, 6).Value = 2003
$ws4.Cells.Item(42, 6).Value = 1337
$ws4.Cells.Item(43, 6).Value = 1209
$ws4.Cells.Item(44, 6).Value = 2108
$ws4.Cells.Item(45, 6).Value = 1131
$ws4.Cells.Item(46, 6).Value = 321
$ws4.Cells.Item(47, 6).Value = 321
$ws4.Cells.Item(48, 6).Value = 2550

$wb.Save()
